$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace row 2 sample data with placeholder values
$ws.Range("B2").Value = "xxxxxx"
$ws.Range("A2").Value = "xxxxx"
$ws.Range("C2").Value = "xxxxx"
$ws.Range("D2").Value = "TZS/ USD"

# Update the selected cell to D2
$ws.Range("D2").Select() | Out-Null
